$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.891.33'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '1.821.08'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9934'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.54%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6283'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9947'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07445'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2926'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.96'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07662'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").Value = '1.827.37'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.972'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6646'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009643'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.026'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '28.893.32'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '224.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9939'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.104'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9942'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '159.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1406'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.461'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.494'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.109'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.042'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05436'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.20%  '
$ws.Range("E32").Value = '  -0.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.847'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7394'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.130'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.609'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").Value = '1.241.94'
$ws.Range("E37").Value = '  -2.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.732'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01770'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.629'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8971'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9941'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.65%  '
$ws.Range("D44").Value = '1.969.43'
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.75'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.74%  '
$ws.Range("E46").Value = '  +2.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5066'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4034'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.918'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.99%  '
$ws.Range("B50").Value = 'XinFinNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07205'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.52%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.655'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.43%  '
